$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "26.000.46" / "210.82" / "1.000" that must stay as
# literal text (European-style thousands separators, significant trailing zeros).
# Mark the whole column as Text before writing so Excel does not reinterpret the
# values as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.000.46"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "1.632.93"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "210.82"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "0.5232"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "0.2584"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").Value = "0.06278"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "20.57"
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D11").Value = "0.07585"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.622.44"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("D13").Value = "4.418"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "1.853.05"
$ws.Range("D15").Value = "0.5495"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "0.0₅8018"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "64.77"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "25.978.28"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "4.676"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").Value = "185.45"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "10.15"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "6.111"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "145.31"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("D27").Value = "7.384"
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").Value = "15.66"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "1.372"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "0.05912"
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("D31").Value = "1.242"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").Value = "3.426"
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("D33").Value = "3.393"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "1.624"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").Value = "0.9807"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.385"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "2.747"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "0.5784"
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("D39").Value = "0.01602"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Value = "0.8487"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "1.038.00"
$ws.Range("E42").Value = "  -5.49%  "
$ws.Range("D43").Value = "5.671"
$ws.Range("E43").Value = "  -7.42%  "
$ws.Range("D44").Value = "100.10"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "1.780.10"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").Value = "54.92"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").Value = "0.9972"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "8.011"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "0.05162"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").Value = "0.4221"
$ws.Range("E51").Value = "  -0.77%  "
